$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update country names that shifted position in the source ranking ---
$ws.Range("A132").Value = 'Lituania'
$ws.Range("A133").Value = 'Gambia'
$ws.Range("A143").Value = 'Georgia'
$ws.Range("A144").Value = 'Estonia'
$ws.Range("A145").Value = 'Sudan del Sur'
$ws.Range("A204").Value = 'Timor Oriental'
$ws.Range("A205").Value = 'Santa Lucia'
$ws.Range("A214").Value = 'Islas Malvinas'
$ws.Range("A215").Value = 'Montserrat'

# --- Update refreshed COVID-19 statistics ---
# Row 28
$ws.Range("B28").Value = 162660
$ws.Range("C28").Value = 2958
$ws.Range("D28").Value = 72324
$ws.Range("E28").Value = 86996
$ws.Range("G28").Value = 76
$ws.Range("H28").Value = 3340
# Row 62
$ws.Range("B62").Value = 46376
$ws.Range("C62").Value = 257
$ws.Range("D62").Value = 42069
$ws.Range("E62").Value = 3384
$ws.Range("G62").Value = 3
$ws.Range("H62").Value = 923
# Row 67
$ws.Range("B67").Value = 38855
$ws.Range("C67").Value = 40
$ws.Range("D67").Value = 32503
$ws.Range("E67").Value = 4916
$ws.Range("G67").Value = 10
$ws.Range("H67").Value = 1436
# Row 88
$ws.Range("B88").Value = 14460
$ws.Range("C88").Value = 581
$ws.Range("D88").Value = 4153
$ws.Range("E88").Value = 9653
$ws.Range("G88").Value = 8
$ws.Range("H88").Value = 654
# Row 132
$ws.Range("B132").Value = 3442
$ws.Range("C132").Value = 45
$ws.Range("D132").Value = 2125
$ws.Range("E132").Value = 1230
$ws.Range("H132").Value = 87
# Row 133
$ws.Range("B133").Value = 3428
$ws.Range("D133").Value = 1737
$ws.Range("E133").Value = 1586
$ws.Range("H133").Value = 105
# Row 143
$ws.Range("B143").Value = 2758
$ws.Range("C143").Value = 196
$ws.Range("D143").Value = 1412
$ws.Range("E143").Value = 1327
$ws.Range("H143").Value = 19
# Row 144
$ws.Range("B144").Value = 2722
$ws.Range("D144").Value = 2286
$ws.Range("E144").Value = 372
$ws.Range("H144").Value = 64
# Row 145
$ws.Range("B145").Value = 2592
$ws.Range("D145").Value = 1290
$ws.Range("E145").Value = 1253
$ws.Range("H145").Value = 49
# Row 160
$ws.Range("B160").Value = 1486
$ws.Range("C160").Value = 4
$ws.Range("E160").Value = 203
# Row 176
$ws.Range("B176").Value = 500
$ws.Range("C176").Value = 1
$ws.Range("D176").Value = 477
# Row 214
$ws.Range("D214").Value = 13
$ws.Range("H214").Value = 0
# Row 215
$ws.Range("D215").Value = 12
$ws.Range("H215").Value = 1

# --- Update "last refreshed" timestamp ---
$ws.Range("A1").Value = 'Datos actualizados a 16 de Septiembre de 2020 a las 09:27'
